$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$asv = $win.ActiveSheetView
$asv | Get-Member | Out-String | Write-Host
